# patient_log.xlsx — row 31 (ID 07071) status update
# The sample has now been reviewed: it previously had no note and a
# "new" status; mark it as not-yet-analyzed and flag its status as
# "error" (matching the other flagged rows, e.g. row 2's "processed").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status (col F) moves from "new" to "error"
$ws.Range("F31").Value = "error"

# Notes (col D) gets a new note explaining the sample hasn't been analyzed
$ws.Range("D31").Value = "haven't analyzed yet"

# Leave the selection on the cell that was just edited
$ws.Range("D31").Select()
